# Update Handback Datetime report values for the "c4ec24ab..." row
# on both the zh-cn and de-de worksheets, as part of regenerating
# the handback status report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-11 08:32:49"
$wsZhCn.Range("H4").Value = "2016-03-11 08:33:36"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-11 08:33:00"
$wsDeDe.Range("H4").Value = "2016-03-11 08:33:55"
